# fix export not contains products bug
#
# The F2 cell held the "${record.product}" placeholder (singular product
# field). The export needs the aggregated product-names string instead, so
# swap the placeholder text for "${record.productNamesString}".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = "`${record.productNamesString}"

# Move/leave the active selection on A3, matching the post-edit cursor spot.
$ws.Range("A3").Select()
